$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("BTC", "Bitcoin", 66967, 1320519938913, 14474144487, -0.17983),
    @("ETH", "Ethereum", 3086.61, 372369245273, 7025145548, -1.28505),
    @("USDT", "Tether", 0.999985, 111432469260, 23784105139, -0.08263),
    @("BNB", "BNB", 578.89, 89096283294, 493976681, 0.02361),
    @("SOL", "Solana", 169.98, 76840909683, 1924984517, -2.20257),
    @("USDC", "USDC", 0.99758, 33479633976, 2282575344, -0.19551),
    @("STETH", "Lido Staked Ether", 3085.81, 28979936067, 55941511, -1.26909),
    @("XRP", "XRP", 0.514706, 28586230718, 380731752, -1.62378),
    @("TON", "Toncoin", 6.42, 22397447049, 200945328, -0.8024),
    @("DOGE", "Dogecoin", 0.15113, 21815222440, 691667227, -2.45723),
    @("ADA", "Cardano", 0.473242, 16752004740, 202247949, -1.45558),
    @("SHIB", "Shiba Inu", 0.00002438, 14408829482, 310501157, -2.14294),
    @("AVAX", "Avalanche", 36.35, 13927846462, 273448625, -2.15134),
    @("TRX", "TRON", 0.121066, 10610782583, 227741888, -2.19123),
    @("WBTC", "Wrapped Bitcoin", 67056, 10445875850, 128505544, -0.12081),
    @("LINK", "Chainlink", 16.51, 9700619238, 490963145, 1.95408),
    @("DOT", "Polkadot", 7.04, 9674917047, 114487709, -1.70956),
    @("BCH", "Bitcoin Cash", 485.47, 9566413219, 320120396, 0.14765),
    @("NEAR", "NEAR Protocol", 7.7, 8338012375, 269472104, -3.45486),
    @("MATIC", "Polygon", 0.691872, 6455330932, 177425773, -3.1812),
    @("LTC", "Litecoin", 83.15000000000001, 6223255909, 186616702, -1.09732),
    @("ICP", "Internet Computer", 12.92, 6035127219, 105609739, -2.83597),
    @("UNI", "Uniswap", 7.74, 5876106544, 214910775, 0.59509),
    @("FET", "Fetch.ai", 2.26, 5732446964, 138116160, -2.74996),
    @("DAI", "Dai", 0.99943, 5505872339, 526143637, -0.04078),
    @("LEO", "LEO Token", 5.88, 5460168590, 766187, 0.36249),
    @("ETC", "Ethereum Classic", 28.03, 4146930109, 122661438, -2.65344),
    @("RNDR", "Render", 10.41, 4072390334, 194213599, 3.89625),
    @("HBAR", "Hedera", 0.111571, 4008554045, 67671445, -1.83067),
    @("PEPE", "Pepe", 0.00000929, 3943199730, 611873523, -6.24193),
    @("APT", "Aptos", 8.380000000000001, 3645497965, 96072083, -2.69209),
    @("WEETH", "Wrapped eETH", 3207.45, 3635545406, 19235699, -1.24172),
    @("FDUSD", "First Digital USD", 0.995981, 3549612286, 2811084669, -0.28767),
    @("IMX", "Immutable", 2.32, 3454320908, 41496758, -3.84599),
    @("CRO", "Cronos", 0.12332, 3300456806, 6037923, -1.9638),
    @("ATOM", "Cosmos Hub", 8.35, 3267454256, 87249589, -3.34919),
    @("FIL", "Filecoin", 5.72, 3189685476, 108343804, -2.8128),
    @("MNT", "Mantle", 0.959202, 3149728739, 67835413, -2.68282),
    @("XLM", "Stellar", 0.106303, 3092143447, 33854070, -2.14642),
    @("EZETH", "Renzo Restaked ETH", 3035.8, 3032838819, 17002935, -1.30834),
    @("AR", "Arweave", 46.29, 3027989864, 55219072, -3.20859),
    @("OKB", "OKB", 49.44, 2976723808, 2699438, -1.4104),
    @("STX", "Stacks", 2.01, 2938141802, 55753791, -3.78026),
    @("KAS", "Kaspa", 0.123, 2920328228, 13372071, -0.1086),
    @("GRT", "The Graph", 0.302826, 2896370076, 58896835, -2.64366),
    @("OP", "Optimism", 2.5, 2736925015, 137954960, -3.12231),
    @("ARB", "Arbitrum", 0.997154, 2655320527, 171042895, -2.34469),
    @("MKR", "Maker", 2776.56, 2590263757, 63695992, -2.45035),
    @("WIF", "dogwifhat", 2.57, 2571125644, 739488697, -4.39164),
    @("TAO", "Bittensor", 377.72, 2565791500, 15533253, -1.00432)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
}
